$p = $ppt.ActivePresentation

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) { return $sh }
    }
    return $null
}

# --- Slide 1 (id 1254): add "Rodriguez" to the author's name ---
$s1 = $p.Slides.Item(1)
$shAuthors = Get-ShapeByName $s1 "Google Shape;1254;p28"
$trAuthors = $shAuthors.TextFrame.TextRange
$fullAuthors = $trAuthors.Text
$nameIdx = $fullAuthors.IndexOf("Bocanegra Hugo Enrique")
$nameChars = $trAuthors.Characters($nameIdx + 1, 22)
$nameChars.Text = "Bocanegra Rodriguez Hugo Enrique"

# --- Slide 2: re-purpose the leftover template slide ---
$s2 = $p.Slides.Item(2)

# Title (id 1340): reposition + retitle ("Componentes importantes a usar",
# first letter kept as its own es-MX run)
$title2 = Get-ShapeByName $s2 "Google Shape;1340;p29"
$title2.Left = 56.30708661417323
$title2.Top = 47.493307186614174
$tr2 = $title2.TextFrame.TextRange
$tr2.Text = "Componentes importantes a usar"
$firstLetter = $tr2.Characters(1, 1)
$firstLetter.Text = "C"
$firstLetter.LanguageID = "es-MX"

# Body placeholder (id 1341): no longer needed, remove entirely
$body2 = Get-ShapeByName $s2 "Google Shape;1341;p29"
$body2.Cut()

# --- Slide 4 (id 1500): flesh out the bullet point ---
$s4 = $p.Slides.Item(4)
$shBullet = Get-ShapeByName $s4 "Google Shape;1500;p32"
$shBullet.TextFrame.TextRange.Text = "- Una idea sobre un prototitpo de casa domotica o automatizada en su mayoria, con la cual ayudara a miles de personas como adultos mayores o personal de trabajo lejos en casa."

# --- Slide 5: update the "goals" section ---
$s5 = $p.Slides.Item(5)

# Accent square (id 1507): nudge into its new position
$accentSquare = Get-ShapeByName $s5 "Google Shape;1507;p33"
$accentSquare.Left = 202.3351212102362
$accentSquare.Top = 123.53393940787402

# Title (id 1508): "OUR GOALS" -> "¿Cómo?"
$title5 = Get-ShapeByName $s5 "Google Shape;1508;p33"
$title5.TextFrame.TextRange.Text = "¿Cómo?"
